# "added request for new TEE Spacer"
#
# Logs a new 3D-print request on the "May 2018" sheet for a TEE Spacer
# (Polyflex material), and tidies up row 27 (the "Date Completed" cell
# had been left un-centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 touch-up: center "Date Completed" like every other cell,
#     and re-pack the row to the sheet's tighter 13.8pt row height.
$ws.Rows(27).RowHeight = 13.8
$ws.Range("B27").HorizontalAlignment = -4108   # xlCenter

# --- New row 29: TEE Spacer print request
$ws.Rows(29).RowHeight = 13.8

$ws.Range("A29").Value = "27-05-2018"
$ws.Range("B29").Value = "27-05-2018"
$ws.Range("C29").Value = "TEE Spacer"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = "Polyflex"
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 10
$ws.Range("H29").Value = 0.2
$ws.Range("I29").Value = "NA"

# --- Leave the selection where the user ended up after logging the row
$ws.Range("I29").Select() | Out-Null
